# "updated records with lab3"
# Fill in LAB THREE (column I) scores on the "Class Quizes" sheet for the
# students who had a grade recorded for that lab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Class Quizes")

$labThreeScores = @{
    2  = 10
    3  = 10
    8  = 15
    9  = 8
    11 = 6
    12 = 10
    13 = 8
    15 = 8
    16 = 10
    17 = 15
    21 = 19
    22 = 15
    24 = 15
    26 = 12
    31 = 19
    32 = 15
    33 = 8
    35 = 15
    38 = 12
    39 = 10
    42 = 8
    43 = 10
    47 = 10
}

foreach ($row in $labThreeScores.Keys) {
    $ws.Cells.Item($row, 9).Value = $labThreeScores[$row]
}

# Leave the workbook focused on the sheet/cell the author last touched.
$null = $ws.Activate()
[void]$ws.Range("I33").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 7

$wb.Application.CalculateFull()
